# Applies the commit:
#  1. Removes the stray empty <c r="B11"> inline-string cell from the
#     "ODI Batting" sheet (row 11, column B).
#  2. Adds a new "ODI Batting Extra" worksheet at the end of the workbook
#     with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# --- 1. ODI Batting: drop the empty B11 cell -----------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B11").ClearContents()

# --- 2. Add the "ODI Batting Extra" sheet at the end ----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
$headerRange = $extraSheet.Range("A1:F1")
for ($col = 1; $col -le 6; $col++) {
    $cell = $extraSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
}
# Reuse the bold / bordered / centered header style already used by the
# other sheets (copy format only, so the shared style is reused rather
# than a near-duplicate being created).
$battingSheet.Range("A1").Copy()
$headerRange.PasteSpecial(-4122)

# Row data: MATCH_CODE, BATTING_POSITION (number or blank), NUM_4 (text),
# NUM_6 (text), PERCENT_RUNS_OF_TOTAL (text), MAN_OF_MATCH (text)
$rows = @(
    @("3045", $null, $null,  $null, $null,    "NO"),
    @("3046", 6,     "0",    "0",   "2.40%",  "NO"),
    @("3462", 7,     "0",    "0",   "2.70%",  "NO"),
    @("3465", 7,     "0",    "0",   "1.35%",  "NO"),
    @("3467", 7,     "0",    "0",   "3.64%",  "NO"),
    @("3738", $null, $null,  $null, $null,    "NO"),
    @("3923", $null, $null,  $null, $null,    "NO"),
    @("3924", 6,     "0",    "0",   "1.94%",  "NO"),
    @("4032", 4,     "4",    "0",   "33.96%", "NO"),
    @("4035", 4,     $null,  $null, $null,    "NO"),
    @("4041", $null, $null,  $null, $null,    "NO"),
    @("4436", $null, $null,  $null, $null,    "NO"),
    @("4437", 4,     "3",    "0",   "7.61%",  "NO"),
    @("4483", 4,     "0",    "0",   "2.78%",  "NO"),
    @("4484", 4,     "0",    "0",   "2.14%",  "NO"),
    @("4486", $null, $null,  $null, $null,    "NO")
)

$r = 2
foreach ($row in $rows) {
    # MATCH_CODE (A) -- numeric-looking, must stay text
    $codeCell = $extraSheet.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    # BATTING_POSITION (B) -- real number when present, else blank text
    $posCell = $extraSheet.Cells.Item($r, 2)
    if ($null -eq $row[1]) {
        $posCell.Value = "'"
    } else {
        $posCell.Value = $row[1]
    }

    # NUM_4 (C), NUM_6 (D), PERCENT_RUNS_OF_TOTAL (E) -- always text
    for ($col = 3; $col -le 5; $col++) {
        $srcVal = $row[$col - 1]
        $cell = $extraSheet.Cells.Item($r, $col)
        if ($null -eq $srcVal) {
            $cell.Value = "'"
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $srcVal
        }
    }

    # MAN_OF_MATCH (F) -- plain text, never numeric-looking
    $extraSheet.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}
